$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.641.69"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").Value = "3.629.56"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'203.22"
$ws.Range("E5").Value = "  +7.80%  "
$ws.Range("D6").Value = "'574.63"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "3.623.78"
$ws.Range("E7").Value = "  +2.44%  "
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "'0.697"
$ws.Range("E10").Value = "  +5.09%  "
$ws.Range("D11").Value = "'61.77"
$ws.Range("E11").Value = "  +17.69%  "
$ws.Range("D12").Value = "'0.152"
$ws.Range("E12").Value = "  +5.01%  "
$ws.Range("D13").Value = "'0.0000288"
$ws.Range("E13").Value = "  +12.49%  "
$ws.Range("D14").Value = "'10.16"
$ws.Range("E14").Value = "  +4.74%  "
$ws.Range("D15").Value = "4.209.24"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "3.629.39"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "'19.14"
$ws.Range("E18").Value = "  +5.48%  "
$ws.Range("D19").Value = "'12.52"
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("D20").Value = "68.469.77"
$ws.Range("E20").Value = "  +3.53%  "
$ws.Range("D21").Value = "'1.08"
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").Value = "'407.43"
$ws.Range("E22").Value = "  +4.47%  "
$ws.Range("D23").Value = "'12.81"
$ws.Range("E23").Value = "  +17.09%  "
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "'86.19"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "'2.95"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'4.01"
$ws.Range("E27").Value = "  +14.90%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'12.69"
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("D29").Value = "'6.16"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").Value = "'9.46"
$ws.Range("E30").Value = "  +7.52%  "
$ws.Range("D31").Value = "'7.94"
$ws.Range("E31").Value = "  +11.65%  "
$ws.Range("D32").Value = "'31.93"
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("D33").Value = "'679.58"
$ws.Range("E33").Value = "  +8.23%  "
$ws.Range("D34").Value = "'12.34"
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D35").Value = "'0.115"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("D36").Value = "'63.83"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").Value = "'42.15"
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("E38").Value = "  +8.62%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0787"
$ws.Range("E40").Value = "  +5.06%  "
$ws.Range("D41").Value = "'3.26"
$ws.Range("E41").Value = "  +17.23%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.213.82"
$ws.Range("E42").Value = "  +8.56%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.136"
$ws.Range("E43").Value = "  +4.53%  "
$ws.Range("D44").Value = "'2.73"
$ws.Range("E44").Value = "  +10.72%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "'2.94"
$ws.Range("E46").Value = "  +27.51%  "
$ws.Range("D47").Value = "'2.92"
$ws.Range("E47").Value = "  +17.19%  "
$ws.Range("E48").Value = "  +5.08%  "
$ws.Range("D49").Value = "'8.93"
$ws.Range("E49").Value = "  +6.69%  "
$ws.Range("D50").Value = "'0.132"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "'3.08"
$ws.Range("E51").Value = "  -1.81%  "
